$d = $word.ActiveDocument

# --- Bootstrap: create numbering.xml (numId=1, decimal list) and the
# "Akapit z listą" / "List Paragraph" paragraph style (styleId "Akapitzlist")
# by driving the same COM calls Word itself uses when a user first presses
# the Numbering button and applies List Paragraph, then tune the generated
# style so it matches a hand-authored "List Paragraph" style exactly
# (basedOn Normalny, uiPriority 34, ind left=720 + contextualSpacing).
$bootstrap = $d.Paragraphs(1)
$bootstrap.Range.Text = "x`r"
$bootstrap.Range.ListFormat.ApplyNumberDefault()
$bootstrap.Range.Style = "Akapit z list"
$listStyle = $d.Styles("Akapit z list")
$listStyle.NameLocal = "List Paragraph"
$listStyle.BaseStyle = "Normalny"
$listStyle.Priority = 34
$listStyle.NoSpaceBetweenParagraphsOfSameStyle = $true
$listStyle.ParagraphFormat.LeftIndent = 36

# --- Replace the whole body (everything before sectPr) with the final
# content in one shot via InsertXML, which lets us control pPr/rPr/run
# layout exactly instead of fighting the Selection/Range text APIs.
$bodyXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:sz w:val="48"/>
          <w:szCs w:val="48"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve">Temat : Wdrożenie systemu </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>informatycznego do szpitala w celu usprawnienia jego procesów.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="Akapitzlist"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Opis dotychczasowego sposobu działania firmy.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="Akapitzlist"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="Akapitzlist"/>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Dotychczasowy sposób działania szpitala obejmował wiele ręcznych i papierowych procesów, co </w:t>
      </w:r>
      <w:r>
        <w:t>prowadziło</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> do błędów, opóźnień i niepotrzebnego obciążenia personelu. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Szpital również </w:t>
      </w:r>
      <w:r>
        <w:t>korzystał z różnych systemów informatycznych, ale nie były one ze sobą zintegrowane, co utrudniało przepływ informacji między różnymi dziedzinami.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="Akapitzlist"/>
        <w:jc w:val="both"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="Akapitzlist"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Opis usprawnień uzyskanych dzięki systemowi informatycznemu.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="Akapitzlist"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="Akapitzlist"/>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:t>Wdrożenie systemu informatycznego w szpitalu może przynieść wiele korzyści. Przede wszystkim, usprawni procesy administracyjne, takie jak rejestracja pacjentów, planowanie wizyt</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> oraz</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> zarządzanie dokumentacją medyczną. System informatyczny może pomóc w automatyzacji tych procesów, co z kolei zminimalizuje ryzyko błędów i opóźnień oraz zmniejszy obciążenie personelu.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>Ponadto, system informatyczny może poprawić bezpieczeństwo</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">. Elektroniczna </w:t>
      </w:r>
      <w:r>
        <w:t>dokumentacj</w:t>
      </w:r>
      <w:r>
        <w:t>a</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> medyczn</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">a oferuje </w:t>
      </w:r>
      <w:r>
        <w:t>łatwiejszy dostęp do informacji o pacjencie, eliminując ryzyko zgubienia lub uszkodzenia</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> co mogło się wydarzyć w przypadku dokumentacji papierowej</w:t>
      </w:r>
      <w:r>
        <w:t>. System informatyczny może również pomóc w lepszym zarządzaniu zasobami</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> m. in.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> sprzęt medyczny, a także w koordynowaniu pracy różnych specjalistów medycznych.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>Reasumując</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">, wdrożenie systemu informatycznego w szpitalu może przynieść wiele korzyści, takich jak zwiększenie efektywności, poprawa jakości opieki zdrowotnej, zmniejszenie kosztów i obciążenia personelu, a także zwiększenie bezpieczeństwa </w:t>
      </w:r>
      <w:r>
        <w:t>dokumentacji medycznej pacjentów.</w:t>
      </w:r>
    </w:p>
'@

$d.Content.InsertXML($bodyXml)
